# Joseph - Training Changes
# Add a new "Changes From 05 May" tracker row documenting:
#   - Developer:           Joseph
#   - Artifacts Name:      Changes to Navigation Files - Service Offering.nav files
#   - How done (SQL col):  Use inserttabconfifuration

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in row 3 of the tracker sheet with the new training change entry.
$ws.Range("E3").Value = "Joseph"
$ws.Range("G3").Value = "Use inserttabconfifuration"
$ws.Range("C3").Value = "Changes to Navigation Files - Service Offering.nav files"

# Move the active selection onto the newly entered cell.
$ws.Range("C3").Select()
